$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append after the existing last row (183 -> row 184/185 in sheet)
$newRows = @(
    @(183, 1, "2024-06-18 20:14:03", 200, 18),
    @(184, 2, "2024-06-18 20:14:04", 200, 3)
)

$startRow = 184
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
